$d = $word.ActiveDocument

# The "Requisitos" section ends with the "LOQ4219: ..." paragraph. In the
# original document it is followed by: a blank paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, the site
# footer/copyright paragraph, and another blank paragraph -- before the
# final (page-break) paragraph that closes the document body. That
# trailing site-footer block (the two text paragraphs plus the blank
# paragraphs immediately around them) is being removed, so the
# "LOQ4219" paragraph becomes immediately followed by the trailing
# page-break paragraph.

$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "LOQ4219:*") {
        $anchorPara = $p
        break
    }
}

$startPara = $anchorPara.Next()

$endPara = $null
$p = $startPara
while ($p -ne $null) {
    if ($p.Format.PageBreakBefore) {
        $endPara = $p
        break
    }
    $p = $p.Next()
}

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
$deleteRange.Delete()
